$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where B (Coin name) and C (Link) change - explicit scalar assignment
$ws.Cells.Item(20, 2).Value = "ProBitToken"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Cells.Item(21, 2).Value = "MCDex"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"

# Data-driven update of Price (D), Volume(1h) (E), Hora (G) columns.
# Each entry: row, new Price, new Volume(1h)%, new Hora
# A leading apostrophe forces Excel to store the value as text (matches the
# workbook's existing inline-string/text cell format for these columns).
$updates = @(
    @(2, "316.62", "3.84%", "8"),
    @(3, "36.17", "0.56%", "8"),
    @(4, "5.151", "1.19%", "8"),
    @(5, "0.08240", "4.85%", "8"),
    @(6, "2.147", "1.44%", "8"),
    @(7, "8.022", "0.70%", "8"),
    @(8, "0.9286", "0.81%", "8"),
    @(9, "0.1012", "4.44%", "8"),
    @(10, "0.1887", "1.34%", "8"),
    @(11, "0.09340", "8.04%", "8"),
    @(12, "0.03633", "3.74%", "8"),
    @(13, "0.09925", "-0.20%", "8"),
    @(14, "0.001439", "0.51%", "8"),
    @(15, "0.005800", "1.37%", "8"),
    @(16, "3.467", "0.12%", "8"),
    @(17, "4.145", "1.14%", "8"),
    @(18, "2.802", "12.35%", "8"),
    @(19, "SKIP", "-1.49%", "8"),
    @(20, "0.1322", "1.59%", "8"),
    @(21, "5.209", "-0.98%", "8"),
    @(22, "0.2252", "2.25%", "8"),
    @(23, "0.04598", "1.25%", "8"),
    @(24, "SKIP", "0.98%", "8"),
    @(25, "0.004732", "-6.93%", "8"),
    @(26, "SKIP", "-21.91%", "8"),
    @(27, "0.0004508", "-5.18%", "8"),
    @(28, "SKIP", "SKIP", "8"),
    @(29, "SKIP", "SKIP", "8"),
    @(30, "SKIP", "SKIP", "8"),
    @(31, "SKIP", "SKIP", "8"),
    @(32, "SKIP", "SKIP", "8"),
    @(33, "SKIP", "SKIP", "8"),
    @(34, "SKIP", "SKIP", "8"),
    @(35, "SKIP", "SKIP", "8"),
    @(36, "SKIP", "SKIP", "8"),
    @(37, "SKIP", "SKIP", "8"),
    @(38, "SKIP", "SKIP", "8"),
    @(39, "0.01998", "7.97%", "8"),
    @(40, "0.04923", "3.87%", "8"),
    @(41, "0.007795", "4.04%", "8"),
    @(42, "0.1401", "-0.05%", "8"),
    @(43, "0.007819", "0.97%", "8"),
    @(44, "0.002106", "-4.52%", "8"),
    @(45, "0.01175", "6.49%", "8"),
    @(46, "0.00006447", "1.07%", "8"),
    @(47, "SKIP", "0.07%", "8"),
    @(48, "40.18", "-2.72%", "8"),
    @(49, "0.001903", "-4.92%", "8"),
    @(50, "SKIP", "0.07%", "8"),
    @(51, "SKIP", "0.07%", "8")
)

foreach ($u in $updates) {
    $row = $u[0]
    if ($u[1] -ne "SKIP") { $ws.Cells.Item($row, 4).Value = "'" + $u[1] }
    if ($u[2] -ne "SKIP") { $ws.Cells.Item($row, 5).Value = "'" + $u[2] }
    if ($u[3] -ne "SKIP") { $ws.Cells.Item($row, 7).Value = "'" + $u[3] }
}